# Auto-generated Excel COM-interop script
# Applies scheduled-runner price/profit refresh to the 8 job-leve sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 33
$ws.Range("H33").Value = 715.2759
$ws.Range("I33").Value = 794.08
$ws.Range("J33").Value = 222.75
$ws.Range("K33").Value = 794.08
$ws.Range("L33").Value = 222.75
$ws.Range("M33").Value = -565.08
$ws.Range("N33").Value = -680.75

# Row 40
$ws.Range("H40").Value = 2491.389
$ws.Range("I40").Value = 2661.8572
$ws.Range("J40").Value = 2252.7334
$ws.Range("K40").Value = 2661.8572
$ws.Range("L40").Value = 2252.7334
$ws.Range("M40").Value = -2486.8572
$ws.Range("N40").Value = -2602.7334

# Row 116
$ws.Range("H116").Value = 1636297.1
$ws.Range("I116").Value = 1986345.9
$ws.Range("J116").Value = 2735.889
$ws.Range("K116").Value = 1986345.9
$ws.Range("L116").Value = 2735.889
$ws.Range("M116").Value = -1982903.9
$ws.Range("N116").Value = -9619.888999999999

# Row 129
$ws.Range("H129").Value = 871.44446
$ws.Range("I129").Value = 320.55554
$ws.Range("J129").Value = 1146.8889
$ws.Range("K129").Value = 961.66662
$ws.Range("L129").Value = 3440.6667
$ws.Range("M129").Value = 4038.33338
$ws.Range("N129").Value = -13440.6667

# Row 132
$ws.Range("H132").Value = 2746.8635
$ws.Range("I132").Value = 2554.4182
$ws.Range("J132").Value = 3709.0908
$ws.Range("K132").Value = 7663.2546
$ws.Range("L132").Value = 11127.2724
$ws.Range("M132").Value = -5133.2546
$ws.Range("N132").Value = -16187.2724

# Row 138
$ws.Range("H138").Value = 2553.8872
$ws.Range("J138").Value = 3033.8655
$ws.Range("L138").Value = 9101.5965
$ws.Range("N138").Value = -19381.5965

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 15539.948
$ws.Range("I32").Value = 4991.2183
$ws.Range("J32").Value = 29353.762
$ws.Range("K32").Value = 4991.2183
$ws.Range("L32").Value = 29353.762
$ws.Range("M32").Value = -4704.2183
$ws.Range("N32").Value = -29927.762

# Row 61
$ws.Range("H61").Value = 1664.7142
$ws.Range("I61").Value = 1664.7142
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 1664.7142
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -1452.7142
$ws.Range("N61").ClearContents()

# Row 94
$ws.Range("H94").Value = 20000
$ws.Range("J94").Value = 20000
$ws.Range("L94").Value = 20000
$ws.Range("N94").Value = -21802

# Row 97
$ws.Range("H97").Value = 2621.9375
$ws.Range("I97").Value = 3449.182
$ws.Range("J97").Value = 802
$ws.Range("K97").Value = 3449.182
$ws.Range("L97").Value = 802
$ws.Range("M97").Value = -2953.182
$ws.Range("N97").Value = -1794

# Row 132
$ws.Range("H132").Value = 2420.3333
$ws.Range("I132").Value = 1777.1177
$ws.Range("J132").Value = 3513.8
$ws.Range("K132").Value = 5331.3531
$ws.Range("L132").Value = 10541.4
$ws.Range("M132").Value = -2801.3531
$ws.Range("N132").Value = -15601.4

# Row 136
$ws.Range("H136").Value = 1664.7142
$ws.Range("I136").Value = 1664.7142
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 4994.142599999999
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -2444.142599999999
$ws.Range("N136").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
# Row 92
$ws.Range("H92").Value = 21950.5
$ws.Range("J92").Value = 21950.5
$ws.Range("L92").Value = 21950.5
$ws.Range("N92").Value = -26942.5

$ws = $wb.Worksheets.Item("CRP")
# Row 106
$ws.Range("H106").Value = 37500
$ws.Range("J106").Value = 37500
$ws.Range("L106").Value = 37500
$ws.Range("N106").Value = -40024

# Row 132
$ws.Range("H132").Value = 1977.8334
$ws.Range("I132").Value = 1262.1177
$ws.Range("J132").Value = 3716
$ws.Range("K132").Value = 3786.3531
$ws.Range("L132").Value = 11148
$ws.Range("M132").Value = -1256.3531
$ws.Range("N132").Value = -16208

# Row 134
$ws.Range("H134").Value = 3418.6487
$ws.Range("I134").Value = 3876.2273
$ws.Range("J134").Value = 2747.5334
$ws.Range("K134").Value = 11628.6819
$ws.Range("L134").Value = 8242.600199999999
$ws.Range("M134").Value = -9093.6819
$ws.Range("N134").Value = -13312.6002

$ws = $wb.Worksheets.Item("CUL")
# Row 122
$ws.Range("H122").Value = 754.94116
$ws.Range("I122").Value = 353
$ws.Range("J122").Value = 1719.6
$ws.Range("K122").Value = 3177
$ws.Range("L122").Value = 15476.4
$ws.Range("M122").Value = -727
$ws.Range("N122").Value = -20376.4

$ws = $wb.Worksheets.Item("GSM")
# Row 2
$ws.Range("H2").Value = 189.71428
$ws.Range("I2").Value = 85.59999999999999
$ws.Range("J2").Value = 450
$ws.Range("K2").Value = 85.59999999999999
$ws.Range("L2").Value = 450
$ws.Range("M2").Value = 27.40000000000001
$ws.Range("N2").Value = -676

# Row 92
$ws.Range("H92").Value = 8011.4165
$ws.Range("J92").Value = 8011.4165
$ws.Range("L92").Value = 8011.4165
$ws.Range("N92").Value = -11755.4165

# Row 95
$ws.Range("H95").Value = 11183.111
$ws.Range("J95").Value = 11183.111
$ws.Range("L95").Value = 11183.111
$ws.Range("N95").Value = -16675.111

# Row 132
$ws.Range("H132").Value = 1718.4375
$ws.Range("I132").Value = 1186.6522
$ws.Range("J132").Value = 3077.4443
$ws.Range("K132").Value = 3559.9566
$ws.Range("L132").Value = 9232.332900000001
$ws.Range("M132").Value = -1029.9566
$ws.Range("N132").Value = -14292.3329

$ws = $wb.Worksheets.Item("LTW")
# Row 104
$ws.Range("H104").Value = 7917.5
$ws.Range("I104").Value = 0
$ws.Range("J104").Value = 7917.5
$ws.Range("K104").Value = 0
$ws.Range("L104").Value = 7917.5
$ws.Range("M104").ClearContents()
$ws.Range("N104").Value = -14905.5

# Row 108
$ws.Range("H108").Value = 64000
$ws.Range("J108").Value = 64000
$ws.Range("L108").Value = 64000
$ws.Range("N108").Value = -71680

$ws = $wb.Worksheets.Item("WVR")
# Row 86
$ws.Range("H86").Value = 16708.334
$ws.Range("I86").Value = 12800
$ws.Range("K86").Value = 12800
$ws.Range("M86").Value = -11677

# Row 89
$ws.Range("H89").Value = 16708.334
$ws.Range("I89").Value = 12800
$ws.Range("K89").Value = 64000
$ws.Range("M89").Value = -58384

# Row 101
$ws.Range("H101").Value = 21602
$ws.Range("J101").Value = 21602
$ws.Range("L101").Value = 21602
$ws.Range("N101").Value = -28092

# Row 104
$ws.Range("H104").Value = 25680
$ws.Range("J104").Value = 25680
$ws.Range("L104").Value = 25680
$ws.Range("N104").Value = -32668

# Row 107
$ws.Range("H107").Value = 488.45
$ws.Range("I107").Value = 210.64285
$ws.Range("K107").Value = 631.9285500000001
$ws.Range("M107").Value = 1288.07145

# Row 132
$ws.Range("H132").Value = 2001.8636
$ws.Range("I132").Value = 628.86664
$ws.Range("J132").Value = 4944
$ws.Range("K132").Value = 1886.59992
$ws.Range("L132").Value = 14832
$ws.Range("M132").Value = 643.4000800000001
$ws.Range("N132").Value = -19892

# Row 136
$ws.Range("H136").Value = 1464.8649
$ws.Range("I136").Value = 976.5769
$ws.Range("J136").Value = 2619
$ws.Range("K136").Value = 2929.7307
$ws.Range("L136").Value = 7857
$ws.Range("M136").Value = -379.7307000000001
$ws.Range("N136").Value = -12957
